{"js": "// Replace the date heading and the division-problem cell texts in document\n// order. Several \"old\" strings repeat (e.g. \"67\u00f72=\" appears twice), and each\n// occurrence maps to a different \"new\" value, so replacements are matched by\n// (text, occurrence-index) rather than a blind global find/replace.\nconst replacements = [\n  [\"2024-05-02 Thursday\", \"2024-05-03 Friday\"],\n  [\"92\u00f74=\", \"36\u00f76=\"],\n  [\"77\u00f75=\", \"36\u00f76=\"],\n  [\"72\u00f72=\", \"19\u00f78=\"],\n  [\"53\u00f77=\", \"33\u00f77=\"],\n  [\"61\u00f72=\", \"25\u00f75=\"],\n  [\"82\u00f74=\", \"89\u00f75=\"],\n  [\"78\u00f77=\", \"23\u00f76=\"],\n  [\"18\u00f75=\", \"52\u00f74=\"],\n  [\"85\u00f77=\", \"55\u00f78=\"],\n  [\"15\u00f79=\", \"84\u00f72=\"],\n  [\"77\u00f79=\", \"16\u00f79=\"],\n  [\"67\u00f72=\", \"61\u00f76=\"],\n  [\"60\u00f74=\", \"26\u00f72=\"],\n  [\"67\u00f72=\", \"72\u00f74=\"],\n  [\"78\u00f72=\", \"16\u00f75=\"],\n  [\"54\u00f72=\", \"91\u00f76=\"],\n  [\"53\u00f72=\", \"79\u00f77=\"],\n  [\"67\u00f75=\", \"93\u00f74=\"],\n  [\"22\u00f76=\", \"33\u00f75=\"],\n  [\"74\u00f78=\", \"90\u00f75=\"],\n  [\"67\u00f73=\", \"44\u00f75=\"],\n  [\"57\u00f77=\", \"40\u00f74=\"],\n  [\"16\u00f74=\", \"81\u00f75=\"],\n  [\"97\u00f75=\", \"22\u00f73=\"],\n  [\"20\u00f74=\", \"27\u00f73=\"],\n];\n\nconst body = context.document.body;\n\n// Collect distinct \"old\" strings, run one search per distinct string, and\n// load all the hits up front.\nconst distinctOld = [...new Set(replacements.map(([oldText]) => oldText))];\nconst searches = {};\nfor (const oldText of distinctOld) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searches[oldText] = found;\n}\n\nawait context.sync();\n\n// Walk the replacement list in document order, consuming one hit per\n// occurrence of a repeated \"old\" string.\nconst nextIndex = {};\nfor (const [oldText, newText] of replacements) {\n  const idx = nextIndex[oldText] || 0;\n  nextIndex[oldText] = idx + 1;\n  const items = searches[oldText].items;\n  items[idx].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-problem cells (5 populated\n# rows x 5 columns in the single table) to the new values. Cell text is set\n# directly by (row, col) position rather than a blind text search, since a\n# couple of the source values (\"67\u00f72=\") repeat but map to different targets\n# depending on which cell they're in.\n$d = $word.ActiveDocument\n\n# Date heading is the document's first paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2024-05-03 Friday\"\n\n$t = $d.Tables.Item(1)\n\n$values = @{\n    1 = @(\"36\u00f76=\", \"36\u00f76=\", \"19\u00f78=\", \"33\u00f77=\", \"25\u00f75=\")\n    5 = @(\"89\u00f75=\", \"23\u00f76=\", \"52\u00f74=\", \"55\u00f78=\", \"84\u00f72=\")\n    9 = @(\"16\u00f79=\", \"61\u00f76=\", \"26\u00f72=\", \"72\u00f74=\", \"16\u00f75=\")\n    13 = @(\"91\u00f76=\", \"79\u00f77=\", \"93\u00f74=\", \"33\u00f75=\", \"90\u00f75=\")\n    17 = @(\"44\u00f75=\", \"40\u00f74=\", \"81\u00f75=\", \"22\u00f73=\", \"27\u00f73=\")\n}\n\nforeach ($row in $values.Keys) {\n    $rowValues = $values[$row]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
